{"js": "// Add a \"Meta description\" paragraph right after the document title, drop\n// the duplicated bold title paragraph near the end of the document, and\n// turn the italic paragraph that used to hold the meta description into\n// the image-generation prompt for the feature image.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1. Insert a new paragraph right after the Heading-1 title with:\n//      [empty run] \"Meta description\" (bold) + \": Looking for a fun, ...\" (regular)\n// ---------------------------------------------------------------------\nconst titlePara = body.paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nconst boldRun = metaPara.insertText(\"Meta description\", \"End\");\nboldRun.font.bold = true;\nawait context.sync();\n\nconst restRun = metaPara.insertText(\n  \": Looking for a fun, free online slot game? Read our review of Action Bank to see why you should play today.\",\n  \"End\"\n);\nrestRun.font.bold = false;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2. Locate (by content, not position) the duplicated bold title\n//    paragraph near the end of the document and remove it, then replace\n//    the text of the italic paragraph that follows it with the new\n//    feature-image prompt (keeping its italic formatting).\n// ---------------------------------------------------------------------\nconst titleHits = body.search(\"Play Action Bank Slot for Free - Read Our Review\", {\n  matchCase: true,\n});\ntitleHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < titleHits.items.length; i++) {\n  const hitPara = titleHits.items[i].paragraphs.getFirst();\n  hitPara.load(\"style\");\n  await context.sync();\n  if (hitPara.style !== \"Heading 1\") {\n    hitPara.delete();\n    await context.sync();\n    break;\n  }\n}\n\nconst descHits = body.search(\n  \"Looking for a fun, free online slot game? Read our review of Action Bank to see why you should play today.\",\n  { matchCase: true }\n);\ndescHits.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < descHits.items.length; i++) {\n  const descPara = descHits.items[i].paragraphs.getFirst();\n  descPara.load(\"text\");\n  await context.sync();\n  // The freshly-inserted meta-description paragraph also contains this\n  // phrase (\"Meta description: Looking for a fun...\"); skip it and only\n  // touch the paragraph whose whole text is the sentence itself.\n  if (descPara.text.trim() === descHits.items[i].text.trim()) {\n    descPara.insertText(\n      'Please design a feature image for the game \"Action Bank\" that fits the following criteria: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses The Maya warrior in the image should be holding a golden vault with a big smile on his face, representing the potential winnings in the game. The warrior should be wearing a traditional Maya headdress and glasses, emphasizing the modern twist to this classic slot game. The background of the image should be bright and colorful, with bold reel symbols including lucky 7s, Xs, bars, and noughts. This feature image should be eye-catching and capture the fun and excitement of playing Action Bank.',\n      \"Replace\"\n    );\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "# Add a \"Meta description\" paragraph right after the document title, drop\n# the duplicated bold title paragraph near the end of the document, and\n# turn the italic paragraph that used to hold the meta description into\n# the image-generation prompt for the feature image.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. Insert a new paragraph right after the Heading-1 title with:\n#      \"Meta description\" (bold) + \": Looking for a fun, ...\" (regular)\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Range.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$metaStart = $metaRange.Start\n\n$boldText = \"Meta description\"\n$metaRange.InsertAfter($boldText)\n$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)\n$boldRange.Bold = 1\n\n$restText = \": Looking for a fun, free online slot game? Read our review of Action Bank to see why you should play today.\"\n$restStart = $metaStart + $boldText.Length\n$restRange = $d.Range($restStart, $restStart)\n$restRange.InsertAfter($restText)\n$restRange2 = $d.Range($restStart, $restStart + $restText.Length)\n$restRange2.Bold = 0\n\n# ---------------------------------------------------------------------\n# 2. Locate (by content, not position) the duplicated bold title\n#    paragraph near the end of the document and remove it, then replace\n#    the text of the italic paragraph that follows it with the new\n#    feature-image prompt (keeping its italic formatting).\n# ---------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -eq \"Play Action Bank Slot for Free - Read Our Review`r\" -and $p.Style.NameLocal -ne \"Heading 1\") {\n    $p.Range.Delete()\n    break\n  }\n}\n\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -eq \"Looking for a fun, free online slot game? Read our review of Action Bank to see why you should play today.`r\") {\n    $newText = 'Please design a feature image for the game \"Action Bank\" that fits the following criteria: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses The Maya warrior in the image should be holding a golden vault with a big smile on his face, representing the potential winnings in the game. The warrior should be wearing a traditional Maya headdress and glasses, emphasizing the modern twist to this classic slot game. The background of the image should be bright and colorful, with bold reel symbols including lucky 7s, Xs, bars, and noughts. This feature image should be eye-catching and capture the fun and excitement of playing Action Bank.'\n    $find = $p.Range.Find\n    $find.Execute($p.Range.Text.TrimEnd(\"`r\"), $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    break\n  }\n}\n"}
